$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename rain model constants (rows 19-24, column A labels).
# Old -> New mapping per commit "Renamed rain model constants":
#   percent_reflecting_sfc -> fi_lidar_rain_reflectivity
#   R                      -> fi_lidar_rain_intensity
#   FaultInjectionLidar4   -> FaultInjectionLidar1
#   FaultInjectionLidar1   -> FaultInjectionLidar2
#   FaultInjectionLidar2   -> FaultInjectionLidar3
#   FaultInjectionLidar3   -> FaultInjectionLidar4
$ws.Range("A20").Value = "fi_lidar_rain_intensity"
$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"
$ws.Range("A21").Value = "FaultInjectionLidar1"
$ws.Range("A22").Value = "FaultInjectionLidar2"
$ws.Range("A23").Value = "FaultInjectionLidar3"
$ws.Range("A24").Value = "FaultInjectionLidar4"

$ws.Range("A20").Select()
